$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 218.75
$ws.Range("I6").Value = 218.75
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 656.25
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -544.25
$ws.Range("N6").ClearContents()

# Row 8
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("K8").Value = 300
$ws.Range("M8").Value = -161

# Row 18
$ws.Range("H18").Value = 1249.75

# Row 33
$ws.Range("H33").Value = 603.1
$ws.Range("I33").Value = 281.22223
$ws.Range("K33").Value = 281.22223
$ws.Range("M33").Value = -52.22223000000002

# Row 38
$ws.Range("H38").Value = 9330.923
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30744

# Row 40
$ws.Range("H40").Value = 2466.3333
$ws.Range("I40").Value = 2200
$ws.Range("J40").Value = 2999
$ws.Range("K40").Value = 2200
$ws.Range("L40").Value = 2999
$ws.Range("M40").Value = -2025
$ws.Range("N40").Value = -3349

# Row 69
$ws.Range("H69").Value = 5015
$ws.Range("J69").Value = 5015
$ws.Range("L69").Value = 15045
$ws.Range("N69").Value = -16793

# Row 72
$ws.Range("H72").Value = 5015
$ws.Range("J72").Value = 5015
$ws.Range("L72").Value = 45135
$ws.Range("N72").Value = -53871

# Row 75
$ws.Range("H75").Value = 94999.5
$ws.Range("J75").Value = 94999.5
$ws.Range("L75").Value = 94999.5
$ws.Range("N75").Value = -96871.5

# Row 78
$ws.Range("H78").Value = 94999.5
$ws.Range("J78").Value = 94999.5
$ws.Range("L78").Value = 284998.5
$ws.Range("N78").Value = -294358.5

# Row 95
$ws.Range("H95").Value = 37666.332
$ws.Range("J95").Value = 37666.332
$ws.Range("L95").Value = 37666.332
$ws.Range("N95").Value = -43158.332

# Row 98
$ws.Range("H98").Value = 1091.6666
$ws.Range("I98").Value = 1091.6666
$ws.Range("K98").Value = 1091.6666
$ws.Range("M98").Value = 406.3334

# Row 122
$ws.Range("H122").Value = 1091.6666
$ws.Range("I122").Value = 1091.6666
$ws.Range("K122").Value = 3274.9998
$ws.Range("M122").Value = -824.9998

# Row 135
$ws.Range("H135").Value = 2118.6
$ws.Range("I135").Value = 1416.6666
$ws.Range("K135").Value = 12749.9994
$ws.Range("M135").Value = -10214.9994

$ws = $wb.Worksheets.Item("ARM")
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 430.8
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 62

# Row 19
$ws.Range("H19").Value = 442.1
$ws.Range("I19").Value = 283.83334
$ws.Range("J19").Value = 679.5
$ws.Range("K19").Value = 283.83334
$ws.Range("L19").Value = 679.5
$ws.Range("M19").Value = -113.83334
$ws.Range("N19").Value = -1019.5

# Row 24
$ws.Range("H24").Value = 442.1
$ws.Range("I24").Value = 283.83334
$ws.Range("J24").Value = 679.5
$ws.Range("K24").Value = 283.83334
$ws.Range("L24").Value = 679.5
$ws.Range("M24").Value = -113.83334
$ws.Range("N24").Value = -1019.5

# Row 31
$ws.Range("H31").Value = 2765.5
$ws.Range("J31").Value = 2500
$ws.Range("L31").Value = 2500
$ws.Range("N31").Value = -3090

# Row 34
$ws.Range("H34").Value = 2765.5
$ws.Range("J34").Value = 2500
$ws.Range("L34").Value = 2500
$ws.Range("N34").Value = -2904

# Row 52
$ws.Range("H52").Value = 85000
$ws.Range("J52").Value = 85000
$ws.Range("L52").Value = 85000
$ws.Range("N52").Value = -85588

# Row 100
$ws.Range("J100").Value = 58994
$ws.Range("L100").Value = 58994
$ws.Range("N100").Value = -61158

# Row 103
$ws.Range("H103").Value = 14000
$ws.Range("I103").Value = 14000
$ws.Range("K103").Value = 14000
$ws.Range("M103").Value = -12828

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 53533.332
$ws.Range("I4").Value = 53533.332
$ws.Range("K4").Value = 160599.996
$ws.Range("M4").Value = -160487.996

# Row 14
$ws.Range("H14").Value = 475
$ws.Range("I14").Value = 475
$ws.Range("K14").Value = 1425
$ws.Range("M14").Value = -1252

# Row 63
$ws.Range("H63").Value = 1833.3334
$ws.Range("J63").Value = 2000
$ws.Range("L63").Value = 6000
$ws.Range("N63").Value = -7498

# Row 66
$ws.Range("H66").Value = 1833.3334
$ws.Range("J66").Value = 2000
$ws.Range("L66").Value = 18000
$ws.Range("N66").Value = -25488

# Row 86
$ws.Range("H86").Value = 412.5
$ws.Range("I86").Value = 400
$ws.Range("J86").Value = 416.66666
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 1249.99998
$ws.Range("M86").Value = -14
$ws.Range("N86").Value = -3621.99998

# Row 89
$ws.Range("H89").Value = 412.5
$ws.Range("I89").Value = 400
$ws.Range("J89").Value = 416.66666
$ws.Range("K89").Value = 3600
$ws.Range("L89").Value = 3749.99994
$ws.Range("M89").Value = 2328
$ws.Range("N89").Value = -15605.99994

# Row 139
$ws.Range("H139").Value = 1256.5
$ws.Range("I139").Value = 1008.6667
$ws.Range("J139").Value = 2000
$ws.Range("K139").Value = 3026.0001
$ws.Range("L139").Value = 6000
$ws.Range("M139").Value = 2113.9999
$ws.Range("N139").Value = -16280

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Row 113
$ws.Range("H113").Value = 5843
$ws.Range("I113").Value = 5843
$ws.Range("K113").Value = 5843
$ws.Range("M113").Value = -3673

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

# Row 100
$ws.Range("H100").Value = 336.5
$ws.Range("I100").Value = 336.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 336.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 204.5
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 4416.5
$ws.Range("I2").Value = 4416.5
$ws.Range("K2").Value = 4416.5
$ws.Range("M2").Value = -4304.5

# Row 105
$ws.Range("H105").Value = 36979.8
$ws.Range("J105").Value = 36979.8
$ws.Range("L105").Value = 36979.8
$ws.Range("N105").Value = -43967.8
